# Updated CHE_grids model - 2025-08-26 18:20
#
# 1) ev_charging_uc!C13 / C14 hold comma-separated timeslice lists that get
#    re-ordered. G7 (=C14) and G8 (=C13) pick the new cached values up via
#    recalculation.
# 2) re_profiles!M4:N7 (season label / weight) rows get shuffled.

$wb = $excel.ActiveWorkbook

$wsEv = $wb.Worksheets.Item("ev_charging_uc")
$wsEv.Range("C13").Value = "FaP,SaP,WaD,RaD,WaP,SaD,RaP,FaD"
$wsEv.Range("C14").Value = "WaP,RaN,FaP,SaP,SaN,WaN,FaN,RaP"

$wsRe = $wb.Worksheets.Item("re_profiles")
$wsRe.Range("M4").Value = "S"
$wsRe.Range("N4").Value = 0.39690767947648675
$wsRe.Range("M5").Value = "F"
$wsRe.Range("N5").Value = 0.27553730672996718
$wsRe.Range("M6").Value = "R"
$wsRe.Range("N6").Value = 0.27149547700006416
$wsRe.Range("M7").Value = "W"
$wsRe.Range("N7").Value = 0.2560595367934817

$wb.Application.Calculate()
